$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("team")

# --- 1. Insert a new row at position 13 (pushes "Juan Fiore" and everyone below down by one) ---
$ws.Rows.Item(13).Insert()

# --- 2 & 3. Populate the new "Juan Douglas" row (13) and finish the "Juan Fiore" row (14) ---
# (cell-write order below matches the original authoring order so new shared strings land
#  at the same indices as the target workbook)
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 3).Value = "Visiting Scholar"

$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 4).Value = "juan_fiore.png"

$ws.Cells.Item(13, 2).Value = "Juan Douglas"
$ws.Cells.Item(13, 4).Value = "juan_douglas.jpg"

$ws.Cells.Item(14, 6).Value = "Agronomist graduated from the National University of Córdoba. Focused on digital agriculture with experience in Python and programming for agricultural applications. Areas of contribution include data management, digital tools, and the integration of technology into agricultural systems."

$ws.Cells.Item(13, 6).Value = "Agronomical Engineering undergraduate student from the national university of Colombia, passionate and skilled in fertigation, irrigation, soil physics, and the use of data in agriculture"

# --- 4. Renumber the id column for all the rows that shifted down ---
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(20, 1).Value = 19

# --- 5. The Hyperlinks collection does not auto-shift with the row insert, so rebuild it ---
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.linkedin.com/in/leonardo-bosche/")
$ws.Range("E3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E2"), "https://www.linkedin.com/in/ignaciociampitti/")
$ws.Range("E2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.linkedin.com/in/natalia-volpato-3824a6150/")
$ws.Range("E5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E9"), "https://www.linkedin.com/in/gmandrini/")
$ws.Range("E9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E18"), "https://www.linkedin.com/in/roberto-carlos-romero-palomeque-831917252?utm_source=share&utm_campaign=share_via&utm_content=profile&utm_medium=ios_app")
$ws.Range("E18").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E15"), "https://www.linkedin.com/in/jjola-unal/", "", "https://www.linkedin.com/in/jjola-unal/")
$ws.Range("E15").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E20"), "https://www.linkedin.com/in/vlasis-m-00537416b")
$ws.Range("E20").Style = "Hyperlink"

# --- 6. Match the saved cursor/selection position from the final workbook ---
$ws.Range("F14").Select()
